# Fruta / hortaliza, semanal
# New weekly price rows for "Feria Lagunitas de Puerto Montt - Kiwi" are
# inserted at the top of this block (row 353), pushing the existing rows
# (old 353-377) down by 3 (new 356-380).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 353 (shifts old rows 353:377 down to 356:380,
# carrying over formatting/number formats from the surrounding rows).
$ws.Range("A353:A355").EntireRow.Insert()

# Values shared by every row in this Kiwi / Feria Lagunitas de Puerto Montt block.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100101007
$categoria   = "Kiwi"
$variedad    = "Hayward"
$unidad      = "$/caja 15 kilos"
$kgUnidad    = 15

function Set-KiwiRow($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg) {
    $ws.Cells.Item($row, 1).Value2 = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value2 = $fecha
    $ws.Cells.Item($row, 5).Value2 = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value2 = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value2 = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value2 = $volumen
    $ws.Cells.Item($row, 14).Value2 = $precioMin
    $ws.Cells.Item($row, 15).Value2 = $precioMax
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value2 = $precioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}

# New row 353: Especial, semana del 2022-09-22 (serial 44826)
Set-KiwiRow 353 44826 "Especial" 100 17000 17000 17000 "Región de O'Higgins" 1133

# New row 354: Primera
Set-KiwiRow 354 44826 "Primera" 100 14000 14000 14000 "Región de O'Higgins" 933

# New row 355: Segunda
Set-KiwiRow 355 44826 "Segunda" 100 12500 12500 12500 "Región de O'Higgins" 833
